{"js": "// Replace the division problems in the document's table with the new values,\n// per the commit diff. Each \"from\" text is unique in the document, so a\n// simple case-sensitive search-and-replace for each pair is safe and\n// order-independent (no pair's replacement text can be accidentally\n// re-matched by a later search, since every \"from\" string is unique).\nconst replacements = [\n  [\"96\u00f77=\", \"45\u00f79=\"],\n  [\"99\u00f77=\", \"69\u00f78=\"],\n  [\"82\u00f72=\", \"95\u00f76=\"],\n  [\"95\u00f74=\", \"31\u00f78=\"],\n  [\"24\u00f79=\", \"11\u00f78=\"],\n  [\"77\u00f79=\", \"34\u00f75=\"],\n  [\"86\u00f79=\", \"34\u00f78=\"],\n  [\"42\u00f72=\", \"74\u00f73=\"],\n  [\"90\u00f77=\", \"70\u00f74=\"],\n  [\"14\u00f78=\", \"71\u00f79=\"],\n  [\"18\u00f76=\", \"57\u00f76=\"],\n  [\"88\u00f72=\", \"56\u00f72=\"],\n  [\"58\u00f72=\", \"45\u00f79=\"],\n  [\"37\u00f78=\", \"41\u00f77=\"],\n  [\"78\u00f76=\", \"80\u00f73=\"],\n  [\"72\u00f79=\", \"80\u00f74=\"],\n  [\"92\u00f74=\", \"86\u00f78=\"],\n  [\"44\u00f79=\", \"12\u00f78=\"],\n  [\"73\u00f72=\", \"38\u00f74=\"],\n  [\"95\u00f79=\", \"26\u00f78=\"],\n  [\"24\u00f75=\", \"94\u00f75=\"],\n  [\"14\u00f76=\", \"27\u00f76=\"],\n  [\"56\u00f74=\", \"89\u00f75=\"],\n  [\"75\u00f74=\", \"91\u00f76=\"],\n  [\"70\u00f75=\", \"51\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division problems in the document's table with the new values,\n# per the commit diff. Each \"from\" text is unique in the document, so a\n# simple Find/Replace (whole-document, case-sensitive, non-wildcard) for each\n# pair is safe and order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ From = \"96\u00f77=\"; To = \"45\u00f79=\" },\n    @{ From = \"99\u00f77=\"; To = \"69\u00f78=\" },\n    @{ From = \"82\u00f72=\"; To = \"95\u00f76=\" },\n    @{ From = \"95\u00f74=\"; To = \"31\u00f78=\" },\n    @{ From = \"24\u00f79=\"; To = \"11\u00f78=\" },\n    @{ From = \"77\u00f79=\"; To = \"34\u00f75=\" },\n    @{ From = \"86\u00f79=\"; To = \"34\u00f78=\" },\n    @{ From = \"42\u00f72=\"; To = \"74\u00f73=\" },\n    @{ From = \"90\u00f77=\"; To = \"70\u00f74=\" },\n    @{ From = \"14\u00f78=\"; To = \"71\u00f79=\" },\n    @{ From = \"18\u00f76=\"; To = \"57\u00f76=\" },\n    @{ From = \"88\u00f72=\"; To = \"56\u00f72=\" },\n    @{ From = \"58\u00f72=\"; To = \"45\u00f79=\" },\n    @{ From = \"37\u00f78=\"; To = \"41\u00f77=\" },\n    @{ From = \"78\u00f76=\"; To = \"80\u00f73=\" },\n    @{ From = \"72\u00f79=\"; To = \"80\u00f74=\" },\n    @{ From = \"92\u00f74=\"; To = \"86\u00f78=\" },\n    @{ From = \"44\u00f79=\"; To = \"12\u00f78=\" },\n    @{ From = \"73\u00f72=\"; To = \"38\u00f74=\" },\n    @{ From = \"95\u00f79=\"; To = \"26\u00f78=\" },\n    @{ From = \"24\u00f75=\"; To = \"94\u00f75=\" },\n    @{ From = \"14\u00f76=\"; To = \"27\u00f76=\" },\n    @{ From = \"56\u00f74=\"; To = \"89\u00f75=\" },\n    @{ From = \"75\u00f74=\"; To = \"91\u00f76=\" },\n    @{ From = \"70\u00f75=\"; To = \"51\u00f77=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.From, $false, $true, $false, $false, $false, $true, 1, $false, $r.To, 2) | Out-Null\n}\n"}
